# Applies cryptos list refresh per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" '29.491.02'
Set-TextValue "D3" '1.607.34'
Set-TextValue "E3" '  +2.84%  '
Set-TextValue "E4" '  +0.00%  '
Set-TextValue "D5" '212.70'
Set-TextValue "E5" '  +1.06%  '
Set-TextValue "D6" '0.522'
Set-TextValue "E6" '  +6.86%  '
Set-TextValue "E7" '  +0.03%  '
Set-TextValue "D8" '26.98'
Set-TextValue "E8" '  +7.35%  '
Set-TextValue "E9" '  -0.74%  '
Set-TextValue "E10" '  +2.68%  '
Set-TextValue "D11" '0.0601'
Set-TextValue "E11" '  +2.55%  '
Set-TextValue "D12" '0.0910'
Set-TextValue "E12" '  +1.56%  '
Set-TextValue "D13" '1.836.02'
Set-TextValue "E13" '  +2.75%  '
Set-TextValue "D14" '1.615.02'
Set-TextValue "E14" '  +3.33%  '
Set-TextValue "D15" '29.505.19'
Set-TextValue "E15" '  +3.14%  '
Set-TextValue "E16" '  +4.05%  '
Set-TextValue "D17" '3.71'
Set-TextValue "E17" '  +1.96%  '
Set-TextValue "D18" '63.23'
Set-TextValue "E18" '  +3.24%  '
Set-TextValue "D19" '241.38'
Set-TextValue "E19" '  +5.20%  '
Set-TextValue "E20" '  +3.74%  '
Set-TextValue "E21" '  +1.79%  '
Set-TextValue "E22" '  +0.02%  '
Set-TextValue "E23" '  +2.58%  '
Set-TextValue "E24" '  +2.33%  '
Set-TextValue "D25" '2.10'
Set-TextValue "E25" '  +0.51%  '
Set-TextValue "D26" '154.71'
Set-TextValue "E26" '  +2.45%  '
Set-TextValue "E27" '  +5.07%  '
Set-TextValue "E28" '  +3.45%  '
Set-TextValue "E29" '  +2.46%  '
Set-TextValue "D30" '0.999'
Set-TextValue "E30" '  -0.01%  '
Set-TextValue "E31" '  +2.58%  '
Set-TextValue "E32" '  +1.17%  '
Set-TextValue "E33" '  +1.76%  '
Set-TextValue "E34" '  +4.17%  '
Set-TextValue "D35" '1.413.66'
Set-TextValue "E35" '  +1.85%  '
Set-TextValue "D36" '1.03'
Set-TextValue "E36" '  +0.60%  '
Set-TextValue "E37" '  +3.38%  '
Set-TextValue "E38" '  +5.00%  '
Set-TextValue "E39" '  +0.34%  '
Set-TextValue "E40" '  +2.55%  '
Set-TextValue "E41" '  +3.60%  '
Set-TextValue "E42" '  +0.84%  '
Set-TextValue "E43" '  +6.01%  '
Set-TextValue "D44" '0.799'
Set-TextValue "E44" '  +3.32%  '
Set-TextValue "B45" 'BitcoinSV'
Set-TextValue "C45" 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue "D45" '52.93'
Set-TextValue "E45" '  +22.43%  '
Set-TextValue "B46" 'PaxDollar'
Set-TextValue "C46" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D46" '0.999'
Set-TextValue "E46" '  +0.00%  '
Set-TextValue "D47" '65.61'
Set-TextValue "E47" '  +2.59%  '
Set-TextValue "D48" '5.29'
Set-TextValue "E48" '  +1.17%  '
Set-TextValue "D49" '1.748.53'
Set-TextValue "E49" '  +2.99%  '
Set-TextValue "D50" '0.856'
Set-TextValue "E50" '  -1.55%  '
Set-TextValue "D51" '86.78'
Set-TextValue "E51" '  +1.84%  '
